$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.286854333333333
$ws.Range("H2").Value = 15.860563
$ws.Range("I2").Value = 0.4075099150374557
$ws.Range("J2").Value = 0.4075099150374557
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 32.017979
$ws.Range("N2").Value = 96.05393700000002
$ws.Range("O2").Value = 0.2161524839374964
$ws.Range("P2").Value = 0.2161524839374964
$ws.Range("Q2").Value = 169.2743910207257
$ws.Range("R2").Value = 1523.469519186531
$ws.Range("S2").Value = 0.08808428036450418
$ws.Range("T2").Value = 0.08808428036450416
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.286854333333333
$ws.Range("H3").Value = 15.860563
$ws.Range("I3").Value = 0.4075099150374557
$ws.Range("J3").Value = 0.4075099150374557
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.90901333333333
$ws.Range("N3").Value = 83.72704
$ws.Range("O3").Value = 0.1884129712323413
$ws.Range("P3").Value = 0.1884129712323412
$ws.Range("Q3").Value = 147.5508880803911
$ws.Range("R3").Value = 1327.95799272352
$ws.Range("S3").Value = 0.07678015389884597
$ws.Range("T3").Value = 0.07678015389884596
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.286854333333333
$ws.Range("H4").Value = 15.860563
$ws.Range("I4").Value = 0.4075099150374557
$ws.Range("J4").Value = 0.4075099150374557
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.29987433333334
$ws.Range("N4").Value = 60.89962300000001
$ws.Range("O4").Value = 0.1370438858982645
$ws.Range("P4").Value = 0.1370438858982645
$ws.Range("Q4").Value = 107.3224785853054
$ws.Range("R4").Value = 965.9023072677491
$ws.Range("S4").Value = 0.05584674229880455
$ws.Range("T4").Value = 0.05584674229880453
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.286854333333333
$ws.Range("H5").Value = 15.860563
$ws.Range("I5").Value = 0.4075099150374557
$ws.Range("J5").Value = 0.4075099150374557
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.468903
$ws.Range("N5").Value = 67.40670900000001
$ws.Range("O5").Value = 0.1516869379794604
$ws.Range("P5").Value = 0.1516869379794604
$ws.Range("Q5").Value = 118.7898171907963
$ws.Range("R5").Value = 1069.108354717167
$ws.Range("S5").Value = 0.06181393120830172
$ws.Range("T5").Value = 0.0618139312083017
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.286854333333333
$ws.Range("H6").Value = 15.860563
$ws.Range("I6").Value = 0.4075099150374557
$ws.Range("J6").Value = 0.4075099150374557
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 17.57385566666667
$ws.Range("N6").Value = 52.721567
$ws.Range("O6").Value = 0.1186406098495176
$ws.Range("P6").Value = 0.1186406098495176
$ws.Range("Q6").Value = 92.91041498469123
$ws.Range("R6").Value = 836.1937348622209
$ws.Range("S6").Value = 0.04834722483976884
$ws.Range("T6").Value = 0.04834722483976883
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.286854333333333
$ws.Range("H7").Value = 15.860563
$ws.Range("I7").Value = 0.4075099150374557
$ws.Range("J7").Value = 0.4075099150374557
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 27.85718966666667
$ws.Range("N7").Value = 83.571569
$ws.Range("O7").Value = 0.1880631111029199
$ws.Range("P7").Value = 0.1880631111029199
$ws.Range("Q7").Value = 147.2769039037052
$ws.Range("R7").Value = 1325.492135133347
$ws.Range("S7").Value = 0.0766375824272305
$ws.Range("T7").Value = 0.07663758242723048
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.035834666666667
$ws.Range("H8").Value = 9.107504
$ws.Range("I8").Value = 0.2340016669801247
$ws.Range("J8").Value = 0.2340016669801247
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 32.017979
$ws.Range("N8").Value = 96.05393700000002
$ws.Range("O8").Value = 0.2161524839374964
$ws.Range("P8").Value = 0.2161524839374964
$ws.Range("Q8").Value = 97.20129060480535
$ws.Range("R8").Value = 874.8116154432482
$ws.Range("S8").Value = 0.0505800415632688
$ws.Range("T8").Value = 0.0505800415632688
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.035834666666667
$ws.Range("H9").Value = 9.107504
$ws.Range("I9").Value = 0.2340016669801247
$ws.Range("J9").Value = 0.2340016669801247
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.90901333333333
$ws.Range("N9").Value = 83.72704
$ws.Range("O9").Value = 0.1884129712323413
$ws.Range("P9").Value = 0.1884129712323412
$ws.Range("Q9").Value = 84.72715018979555
$ws.Range("R9").Value = 762.5443517081601
$ws.Range("S9").Value = 0.04408894934904614
$ws.Range("T9").Value = 0.04408894934904613
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.035834666666667
$ws.Range("H10").Value = 9.107504
$ws.Range("I10").Value = 0.2340016669801247
$ws.Range("J10").Value = 0.2340016669801247
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.29987433333334
$ws.Range("N10").Value = 60.89962300000001
$ws.Range("O10").Value = 0.1370438858982645
$ws.Range("P10").Value = 0.1370438858982645
$ws.Range("Q10").Value = 61.62706223011023
$ws.Range("R10").Value = 554.6435600709921
$ws.Range("S10").Value = 0.0320684977496279
$ws.Range("T10").Value = 0.0320684977496279
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.035834666666667
$ws.Range("H11").Value = 9.107504
$ws.Range("I11").Value = 0.2340016669801247
$ws.Range("J11").Value = 0.2340016669801247
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 22.468903
$ws.Range("N11").Value = 67.40670900000001
$ws.Range("O11").Value = 0.1516869379794604
$ws.Range("P11").Value = 0.1516869379794604
$ws.Range("Q11").Value = 68.21187464937067
$ws.Range("R11").Value = 613.9068718443361
$ws.Range("S11").Value = 0.03549499634630453
$ws.Range("T11").Value = 0.03549499634630452
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.035834666666667
$ws.Range("H12").Value = 9.107504
$ws.Range("I12").Value = 0.2340016669801247
$ws.Range("J12").Value = 0.2340016669801247
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 17.57385566666667
$ws.Range("N12").Value = 52.721567
$ws.Range("O12").Value = 0.1186406098495176
$ws.Range("P12").Value = 0.1186406098495176
$ws.Range("Q12").Value = 53.35132025986311
$ws.Range("R12").Value = 480.161882338768
$ws.Range("S12").Value = 0.02776210047632572
$ws.Range("T12").Value = 0.02776210047632571
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.035834666666667
$ws.Range("H13").Value = 9.107504
$ws.Range("I13").Value = 0.2340016669801247
$ws.Range("J13").Value = 0.2340016669801247
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 27.85718966666667
$ws.Range("N13").Value = 83.571569
$ws.Range("O13").Value = 0.1880631111029199
$ws.Range("P13").Value = 0.1880631111029199
$ws.Range("Q13").Value = 84.56982210597511
$ws.Range("R13").Value = 761.128398953776
$ws.Range("S13").Value = 0.04400708149555167
$ws.Range("T13").Value = 0.04400708149555167
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 4.650871
$ws.Range("H14").Value = 13.952613
$ws.Range("I14").Value = 0.3584884179824196
$ws.Range("J14").Value = 0.3584884179824196
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 32.017979
$ws.Range("N14").Value = 96.05393700000002
$ws.Range("O14").Value = 0.2161524839374964
$ws.Range("P14").Value = 0.2161524839374964
$ws.Range("Q14").Value = 148.911490009709
$ws.Range("R14").Value = 1340.203410087381
$ws.Range("S14").Value = 0.07748816200972346
$ws.Range("T14").Value = 0.07748816200972347
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 4.650871
$ws.Range("H15").Value = 13.952613
$ws.Range("I15").Value = 0.3584884179824196
$ws.Range("J15").Value = 0.3584884179824196
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 27.90901333333333
$ws.Range("N15").Value = 83.72704
$ws.Range("O15").Value = 0.1884129712323413
$ws.Range("P15").Value = 0.1884129712323412
$ws.Range("Q15").Value = 129.8012207506133
$ws.Range("R15").Value = 1168.21098675552
$ws.Range("S15").Value = 0.06754386798444915
$ws.Range("T15").Value = 0.06754386798444915
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 4.650871
$ws.Range("H16").Value = 13.952613
$ws.Range("I16").Value = 0.3584884179824196
$ws.Range("J16").Value = 0.3584884179824196
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 20.29987433333334
$ws.Range("N16").Value = 60.89962300000001
$ws.Range("O16").Value = 0.1370438858982645
$ws.Range("P16").Value = 0.1370438858982645
$ws.Range("Q16").Value = 94.41209684054434
$ws.Range("R16").Value = 849.708871564899
$ws.Range("S16").Value = 0.04912864584983207
$ws.Range("T16").Value = 0.04912864584983207
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 4.650871
$ws.Range("H17").Value = 13.952613
$ws.Range("I17").Value = 0.3584884179824196
$ws.Range("J17").Value = 0.3584884179824196
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.468903
$ws.Range("N17").Value = 67.40670900000001
$ws.Range("O17").Value = 0.1516869379794604
$ws.Range("P17").Value = 0.1516869379794604
$ws.Range("Q17").Value = 104.499969364513
$ws.Range("R17").Value = 940.4997242806171
$ws.Range("S17").Value = 0.05437801042485416
$ws.Range("T17").Value = 0.05437801042485416
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 4.650871
$ws.Range("H18").Value = 13.952613
$ws.Range("I18").Value = 0.3584884179824196
$ws.Range("J18").Value = 0.3584884179824196
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 17.57385566666667
$ws.Range("N18").Value = 52.721567
$ws.Range("O18").Value = 0.1186406098495176
$ws.Range("P18").Value = 0.1186406098495176
$ws.Range("Q18").Value = 81.73373567828565
$ws.Range("R18").Value = 735.6036211045709
$ws.Range("S18").Value = 0.04253128453342302
$ws.Range("T18").Value = 0.04253128453342302
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 4.650871
$ws.Range("H19").Value = 13.952613
$ws.Range("I19").Value = 0.3584884179824196
$ws.Range("J19").Value = 0.3584884179824196
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 27.85718966666667
$ws.Range("N19").Value = 83.571569
$ws.Range("O19").Value = 0.1880631111029199
$ws.Range("P19").Value = 0.1880631111029199
$ws.Range("Q19").Value = 129.5601955621997
$ws.Range("R19").Value = 1166.041760059797
$ws.Range("S19").Value = 0.06741844718013777
$ws.Range("T19").Value = 0.06741844718013779
